$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1032
$ws.Range("I2").Value = 421.5
$ws.Range("K2").Value = 421.5
$ws.Range("M2").Value = -308.5
$ws.Range("H86").Value = 2364.2354
$ws.Range("J86").Value = 2410.1428
$ws.Range("L86").Value = 2410.1428
$ws.Range("N86").Value = -4656.1428
$ws.Range("H89").Value = 2364.2354
$ws.Range("J89").Value = 2410.1428
$ws.Range("L89").Value = 12050.714
$ws.Range("N89").Value = -23282.714
$ws.Range("H106").Value = 2082.5
$ws.Range("I106").Value = 1808.5714
$ws.Range("K106").Value = 1808.5714
$ws.Range("M106").Value = -1177.5714
$ws.Range("H116").Value = 5517.5
$ws.Range("I116").Value = 5765.7144
$ws.Range("J116").Value = 5269.2856
$ws.Range("K116").Value = 5765.7144
$ws.Range("L116").Value = 5269.2856
$ws.Range("M116").Value = -2323.7144
$ws.Range("N116").Value = -12153.2856
$ws.Range("H141").Value = 1058.3334
$ws.Range("I141").Value = 1058.3334
$ws.Range("K141").Value = 3175.0002
$ws.Range("M141").Value = 2004.9998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I74").Value = 3002006
$ws.Range("J74").Value = 6249.5
$ws.Range("K74").Value = 3002006
$ws.Range("L74").Value = 6249.5
$ws.Range("M74").Value = -3001132
$ws.Range("N74").Value = -7997.5
$ws.Range("I77").Value = 3002006
$ws.Range("J77").Value = 6249.5
$ws.Range("K77").Value = 15010030
$ws.Range("L77").Value = 31247.5
$ws.Range("M77").Value = -15005662
$ws.Range("N77").Value = -39983.5
$ws.Range("H109").Value = 103000
$ws.Range("J109").Value = 103000
$ws.Range("L109").Value = 103000
$ws.Range("N109").Value = -105774
$ws.Range("H132").Value = 1604.2593
$ws.Range("I132").Value = 832.6818
$ws.Range("J132").Value = 4999.2
$ws.Range("K132").Value = 2498.0454
$ws.Range("L132").Value = 14997.6
$ws.Range("M132").Value = 31.95460000000003
$ws.Range("N132").Value = -20057.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 732.7778
$ws.Range("J94").Value = 610.6667
$ws.Range("L94").Value = 610.6667
$ws.Range("N94").Value = -1512.6667
$ws.Range("H99").Value = 2559.5386
$ws.Range("I99").Value = 1841.2222
$ws.Range("J99").Value = 4175.75
$ws.Range("K99").Value = 1841.2222
$ws.Range("L99").Value = 4175.75
$ws.Range("M99").Value = -343.2221999999999
$ws.Range("N99").Value = -7171.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 166.21428
$ws.Range("I7").Value = 146.09091
$ws.Range("J7").Value = 240
$ws.Range("K7").Value = 146.09091
$ws.Range("L7").Value = 240
$ws.Range("M7").Value = -33.09091000000001
$ws.Range("N7").Value = -466
$ws.Range("H31").Value = 2704366.2
$ws.Range("I31").Value = 3449527.8
$ws.Range("J31").Value = 3156.125
$ws.Range("K31").Value = 3449527.8
$ws.Range("L31").Value = 3156.125
$ws.Range("M31").Value = -3449232.8
$ws.Range("N31").Value = -3746.125
$ws.Range("H34").Value = 2704366.2
$ws.Range("I34").Value = 3449527.8
$ws.Range("J34").Value = 3156.125
$ws.Range("K34").Value = 3449527.8
$ws.Range("L34").Value = 3156.125
$ws.Range("M34").Value = -3449325.8
$ws.Range("N34").Value = -3560.125
$ws.Range("H38").Value = 3233
$ws.Range("J38").Value = 3233
$ws.Range("L38").Value = 3233
$ws.Range("N38").Value = -3987
$ws.Range("H41").Value = 23744.5
$ws.Range("I41").Value = 1211.2
$ws.Range("J41").Value = 36263
$ws.Range("K41").Value = 1211.2
$ws.Range("L41").Value = 36263
$ws.Range("M41").Value = -783.2
$ws.Range("N41").Value = -37119
$ws.Range("H46").Value = 3233
$ws.Range("J46").Value = 3233
$ws.Range("L46").Value = 3233
$ws.Range("N46").Value = -3655
$ws.Range("H50").Value = 34498.223
$ws.Range("J50").Value = 34498.223
$ws.Range("L50").Value = 34498.223
$ws.Range("N50").Value = -35748.223
$ws.Range("H51").Value = 36979.2
$ws.Range("J51").Value = 36979.2
$ws.Range("L51").Value = 36979.2
$ws.Range("N51").Value = -38451.2
$ws.Range("H58").Value = 26648
$ws.Range("I58").Value = 2011
$ws.Range("J58").Value = 66067.2
$ws.Range("K58").Value = 2011
$ws.Range("L58").Value = 66067.2
$ws.Range("M58").Value = -1808
$ws.Range("N58").Value = -66473.2
$ws.Range("H59").Value = 26574.8
$ws.Range("J59").Value = 26969
$ws.Range("L59").Value = 26969
$ws.Range("N59").Value = -29259
$ws.Range("H60").Value = 33288.168
$ws.Range("J60").Value = 33046
$ws.Range("L60").Value = 33046
$ws.Range("N60").Value = -34068
$ws.Range("H61").Value = 36979.2
$ws.Range("J61").Value = 36979.2
$ws.Range("L61").Value = 36979.2
$ws.Range("N61").Value = -37675.2
$ws.Range("H120").Value = 19486.2
$ws.Range("I120").Value = 19431
$ws.Range("K120").Value = 19431
$ws.Range("M120").Value = -15802
$ws.Range("H122").Value = 1132.2
$ws.Range("I122").Value = 977.5714
$ws.Range("K122").Value = 2932.7142
$ws.Range("M122").Value = -482.7142000000003
$ws.Range("H134").Value = 1802.5532
$ws.Range("I134").Value = 1554.9
$ws.Range("K134").Value = 4664.700000000001
$ws.Range("M134").Value = -2129.700000000001
$ws.Range("H136").Value = 26648
$ws.Range("I136").Value = 2011
$ws.Range("J136").Value = 66067.2
$ws.Range("K136").Value = 6033
$ws.Range("L136").Value = 198201.6
$ws.Range("M136").Value = -3483
$ws.Range("N136").Value = -203301.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 803.4545000000001
$ws.Range("I5").Value = 783.8
$ws.Range("K5").Value = 2351.4
$ws.Range("M5").Value = -2239.4
$ws.Range("H6").Value = 7545.75
$ws.Range("I6").Value = 6729.6665
$ws.Range("K6").Value = 20188.9995
$ws.Range("M6").Value = -20075.9995
$ws.Range("H135").Value = 803.4545000000001
$ws.Range("I135").Value = 783.8
$ws.Range("K135").Value = 7054.2
$ws.Range("M135").Value = -4519.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 50000
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H46").Value = 33999.5
$ws.Range("J46").Value = 37999
$ws.Range("L46").Value = 37999
$ws.Range("N46").Value = -38311

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1320.2106
$ws.Range("I22").Value = 1334.909
$ws.Range("K22").Value = 1334.909
$ws.Range("M22").Value = -1039.909
$ws.Range("H27").Value = 1320.2106
$ws.Range("I27").Value = 1334.909
$ws.Range("K27").Value = 1334.909
$ws.Range("M27").Value = -1227.909
$ws.Range("H61").Value = 631.2
$ws.Range("I61").Value = 623
$ws.Range("K61").Value = 623
$ws.Range("M61").Value = -421
$ws.Range("H113").Value = 631.2
$ws.Range("I113").Value = 623
$ws.Range("K113").Value = 623
$ws.Range("M113").Value = 1547
$ws.Range("H136").Value = 4758.294
$ws.Range("I136").Value = 4739.0347
$ws.Range("K136").Value = 14217.1041
$ws.Range("M136").Value = -11667.1041

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 21711.5
$ws.Range("I54").Value = 6000
$ws.Range("K54").Value = 6000
$ws.Range("M54").Value = -5480
$ws.Range("H63").Value = 49233.75
$ws.Range("J63").Value = 49233.75
$ws.Range("L63").Value = 49233.75
$ws.Range("N63").Value = -50481.75
$ws.Range("H66").Value = 49233.75
$ws.Range("J66").Value = 49233.75
$ws.Range("L66").Value = 147701.25
$ws.Range("N66").Value = -153941.25
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H109").Value = 57492
$ws.Range("I109").Value = 19999
$ws.Range("J109").Value = 64990.6
$ws.Range("K109").Value = 19999
$ws.Range("L109").Value = 64990.6
$ws.Range("M109").Value = -18612
$ws.Range("N109").Value = -67764.60000000001
$ws.Range("H113").Value = 1274.6666
$ws.Range("I113").Value = 998.1667
$ws.Range("J113").Value = 1827.6666
$ws.Range("K113").Value = 2994.5001
$ws.Range("L113").Value = 5482.9998
$ws.Range("M113").Value = -824.5001000000002
$ws.Range("N113").Value = -9822.9998
$ws.Range("H136").Value = 17658.871
$ws.Range("I136").Value = 18343.217
$ws.Range("J136").Value = 4998.5
$ws.Range("K136").Value = 55029.651
$ws.Range("L136").Value = 14995.5
$ws.Range("M136").Value = -52479.651
$ws.Range("N136").Value = -20095.5
